# ranking/historico.xlsx — "Add files via upload"
#
# The sheet (Planilha1, A1:E51) is a leaderboard with columns:
#   A = id, B = name, C = inicial, D = posicao, E = link
#
# This edit:
#   1. Renames the "name" of the row-22 entry from "Rafael Andrés" to
#      "MDR Player" (same id/link, same rank).
#   2. Inserts a brand-new ranked entry as the new row 50 (posicao 49),
#      pushing the two former tail rows (49->50, 50->51 posicao) down by
#      one row each, to position 51/52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 22: rename the player from "Rafael Andrés" to "MDR Player".
$ws.Range("B22").Value2 = "MDR Player"

# 2. Insert a fresh row at 50, shifting the old rows 50 & 51 down to 51 & 52.
$ws.Rows("50:50").Insert(-4121)          # -4121 = xlShiftDown
$ws.Range("A50:E50").ClearFormats()      # new row carries no cell style

$ws.Range("A50").Value2 = "66712ed43d409a28fd524ff3"
$ws.Range("B50").Value2 = "Rafael Andrés"
$ws.Range("C50").Value2 = 6053416.017
$ws.Range("D50").Value2 = 49
$ws.Range("E50").Value2 = "RFaNFT"

# The "posicao" column is static data, not a formula, so fix up the two
# rows that got pushed down a slot.
$ws.Range("D51").Value2 = 50
$ws.Range("D52").Value2 = 51

# Keep the view/selection consistent with the new, larger range.
$null = $ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("A2:E52").Select()

Write-Host "done"
